$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at A. This shifts every existing column one to
# the right (old A -> B, old B -> C, ... old Q -> R) and carries the
# existing cell values/styles along with it.
$ws.Columns("A:A").Insert()

# The old column A ("Status" header / {bookings:internal_status} placeholder)
# is now column B. Rename its header from "Status" to "Internal Status".
# (Set this before writing the brand-new strings below so the shared
# string table ends up in the same append order as the authored edit.)
$ws.Range("B1").Value = "Internal Status"

# Populate the new column A with the "Current Status" header/placeholder.
$ws.Range("A1").Value = "Current Status"
$ws.Range("A2").Value = "{bookings:current_status}"

# Match formatting of the other header cells (bold) for the two new/renamed
# header cells.
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Name = "Arial"
$ws.Range("A1:B1").Font.Size = 10

# Placeholder row formatting: keep it consistent with its neighbouring
# placeholder cells (non-bold Arial 10); B2 (Internal Status placeholder)
# is centered like the other placeholder cells in that row.
$ws.Range("A2:B2").Font.Name = "Arial"
$ws.Range("A2:B2").Font.Size = 10
$ws.Range("B2").HorizontalAlignment = -4108

# Column widths: new column A ("Current Status") and column B ("Internal
# Status") get their own explicit widths; column C ("Booking Date") also
# widens slightly to fit the new layout.
$ws.Columns("A").ColumnWidth = 21.833333333333332
$ws.Columns("B").ColumnWidth = 27.333333333333332
$ws.Columns("C").ColumnWidth = 20.833333333333332

# Move the active selection to A2 (matches the saved view state).
$ws.Range("A2").Select()
